# Refresh the "Data" sheet's 12-row result block (rows 2-13) with the
# latest batch of lab-run metadata: 20220915-Salm-13035Updt /
# CartridgeSalm3035 / Result IDs A1528001-A1528012.
#
# Result ID (col A), Lab Sample ID (col E) and Cartridge ID (col T) are the
# only fields that change; everything else in the row stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$batchName     = "20220915-Salm-13035Updt"
$cartridgeName = "CartridgeSalm3035"

$resultIds = @(
    "A1528001", "A1528002", "A1528003", "A1528004",
    "A1528005", "A1528006", "A1528007", "A1528008",
    "A1528009", "A1528010", "A1528011", "A1528012"
)

for ($i = 0; $i -lt $resultIds.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value  = $resultIds[$i]   # A: Result ID
    $ws.Cells.Item($row, 5).Value  = $batchName        # E: Lab Sample ID
    $ws.Cells.Item($row, 20).Value = $cartridgeName     # T: Cartridge ID
}
